$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "29.822.06"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "1.885.17"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7476"
$ws.Range("E5").Value = "  -3.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.89"
$ws.Range("E6").Value = "  -0.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3120"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.16"
$ws.Range("E9").Value = "  -2.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07095"
$ws.Range("E10").Value = "  -3.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08496"
$ws.Range("E11").Value = "  +5.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7576"
$ws.Range("E12").Value = "  -1.82%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.351"
$ws.Range("E13").Value = "  -2.66%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.822.08"
$ws.Range("E14").Value = "  -5.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.88"
$ws.Range("E15").Value = "  -1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.118"
$ws.Range("E16").Value = "  -1.70%  "

$ws.Range("D17").Value = "29.781.08"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.67"
$ws.Range("E18").Value = "  -2.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.53"
$ws.Range("E19").Value = "  -1.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007819"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9991"
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").Value = "2.139.68"
$ws.Range("E22").Value = "  -0.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.964"
$ws.Range("E23").Value = "  -2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1585"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.337"
$ws.Range("E26").Value = "  -1.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.86"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.67"
$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.026"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("E30").Value = "  +3.23%  "

$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.493"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.158"
$ws.Range("E33").Value = "  +2.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05406"
$ws.Range("E34").Value = "  -3.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.236"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7507"
$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.706"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01941"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.768"
$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4455"
$ws.Range("E41").Value = "  -0.32%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.082"
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.096.75"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.07"
$ws.Range("E44").Value = "  -3.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8595"
$ws.Range("E45").Value = "  +1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.707"
$ws.Range("E47").Value = "  +2.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.31"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.853"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.028"
$ws.Range("E50").Value = "  +0.84%  "

$ws.Range("D51").Value = "2.036.93"
$ws.Range("E51").Value = "  -1.14%  "
